$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing ImagesCount values (B2:B5) from 4 to 2 ---
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

# --- Fill in species/counts for rows 6-9 (cells already exist, keep row style) ---
$ws.Range("A6").Value = "w3"
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = "w3"
$ws.Range("B7").Value = 2
$ws.Range("B7").Style = "Normal"

$ws.Range("A8").Value = "w4"
$ws.Range("B8").Value = 2

$ws.Range("A9").Value = "w4"
$ws.Range("B9").Value = 2

# --- New data rows 10-19 (previously empty rows, cells get default/no style) ---
$ws.Range("A10").Value = "w5"
$ws.Range("B10").Value = 2
$ws.Range("A10:B10").Style = "Normal"

$ws.Range("A11").Value = "w5"
$ws.Range("B11").Value = 2
$ws.Range("A11:B11").Style = "Normal"

$ws.Range("A12").Value = "w6"
$ws.Range("B12").Value = 2
$ws.Range("A12:B12").Style = "Normal"

$ws.Range("A13").Value = "w6"
$ws.Range("B13").Value = 2
$ws.Range("A13:B13").Style = "Normal"

$ws.Range("A14").Value = "w7"
$ws.Range("B14").Value = 2
$ws.Range("A14:B14").Style = "Normal"

$ws.Range("A15").Value = "w7"
$ws.Range("B15").Value = 2
$ws.Range("A15:B15").Style = "Normal"

$ws.Range("A16").Value = "w8"
$ws.Range("B16").Value = 2
$ws.Range("A16:B16").Style = "Normal"

$ws.Range("A17").Value = "w8"
$ws.Range("B17").Value = 2
$ws.Range("A17:B17").Style = "Normal"

$ws.Range("A18").Value = "w9"
$ws.Range("B18").Value = 2
$ws.Range("A18:B18").Style = "Normal"

$ws.Range("A19").Value = "w9"
$ws.Range("B19").Value = 2
$ws.Range("A19:B19").Style = "Normal"

# --- Update the active selection on the sheet ---
$null = $ws.Range("A2:B5").Select()

# --- Update workbook window geometry (best-effort; runtime may ignore) ---
$win = $wb.Windows.Item(1)
$win.Left = 1480
$win.Top = 1480
$win.Width = 14400
$win.Height = 7360
